$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gdnf"
$ws.Range("C2").Value = "Ret"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1170883333333333
$ws.Range("H2").Value = 0.351265
$ws.Range("I2").Value = 0.09352043543917719
$ws.Range("J2").Value = 0.09352043543917718
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.273058
$ws.Range("N2").Value = 3.819174
$ws.Range("O2").Value = 0.1826308343983922
$ws.Range("P2").Value = 0.1826308343983922
$ws.Range("Q2").Value = 0.1490602394566667
$ws.Range("R2").Value = 1.34154215511
$ws.Range("S2").Value = 0.0170797151575579
$ws.Range("T2").Value = 0.0170797151575579

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gdnf"
$ws.Range("C3").Value = "Ret"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1170883333333333
$ws.Range("H3").Value = 0.351265
$ws.Range("I3").Value = 0.09352043543917719
$ws.Range("J3").Value = 0.09352043543917718
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.094146333333334
$ws.Range("N3").Value = 15.282439
$ws.Range("O3").Value = 0.7307979647464429
$ws.Range("P3").Value = 0.7307979647464428
$ws.Range("Q3").Value = 0.5964651039261112
$ws.Range("R3").Value = 5.368185935335
$ws.Range("S3").Value = 0.0683445438811518
$ws.Range("T3").Value = 0.06834454388115178

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gdnf"
$ws.Range("C4").Value = "Ret"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1170883333333333
$ws.Range("H4").Value = 0.351265
$ws.Range("I4").Value = 0.09352043543917719
$ws.Range("J4").Value = 0.09352043543917718
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.012811
$ws.Range("N4").Value = 0.038433
$ws.Range("O4").Value = 0.00183784526665541
$ws.Range("P4").Value = 0.00183784526665541
$ws.Range("Q4").Value = 0.001500018638333334
$ws.Range("R4").Value = 0.013500167745
$ws.Range("S4").Value = 0.0001718760896074447
$ws.Range("T4").Value = 0.0001718760896074446

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Gdnf"
$ws.Range("C5").Value = "Ret"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1170883333333333
$ws.Range("H5").Value = 0.351265
$ws.Range("I5").Value = 0.09352043543917719
$ws.Range("J5").Value = 0.09352043543917718
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5906476666666668
$ws.Range("N5").Value = 1.771943
$ws.Range("O5").Value = 0.08473335558850956
$ws.Range("P5").Value = 0.08473335558850954
$ws.Range("Q5").Value = 0.06915795087722224
$ws.Range("R5").Value = 0.6224215578950001
$ws.Range("S5").Value = 0.007924300310860052
$ws.Range("T5").Value = 0.00792430031086005

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Gdnf"
$ws.Range("C6").Value = "Ret"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.134919666666667
$ws.Range("H6").Value = 3.404758999999999
$ws.Range("I6").Value = 0.9064795645608229
$ws.Range("J6").Value = 0.9064795645608228
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.273058
$ws.Range("N6").Value = 3.819174
$ws.Range("O6").Value = 0.1826308343983922
$ws.Range("P6").Value = 0.1826308343983922
$ws.Range("Q6").Value = 1.444818561007333
$ws.Range("R6").Value = 13.003367049066
$ws.Range("S6").Value = 0.1655511192408343
$ws.Range("T6").Value = 0.1655511192408343

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Gdnf"
$ws.Range("C7").Value = "Ret"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.134919666666667
$ws.Range("H7").Value = 3.404758999999999
$ws.Range("I7").Value = 0.9064795645608229
$ws.Range("J7").Value = 0.9064795645608228
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.094146333333334
$ws.Range("N7").Value = 15.282439
$ws.Range("O7").Value = 0.7307979647464429
$ws.Range("P7").Value = 0.7307979647464428
$ws.Range("Q7").Value = 5.781446858577889
$ws.Range("R7").Value = 52.03302172720099
$ws.Range("S7").Value = 0.6624534208652911
$ws.Range("T7").Value = 0.6624534208652909

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gdnf"
$ws.Range("C8").Value = "Ret"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.134919666666667
$ws.Range("H8").Value = 3.404758999999999
$ws.Range("I8").Value = 0.9064795645608229
$ws.Range("J8").Value = 0.9064795645608228
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.012811
$ws.Range("N8").Value = 0.038433
$ws.Range("O8").Value = 0.00183784526665541
$ws.Range("P8").Value = 0.00183784526665541
$ws.Range("Q8").Value = 0.01453945584966667
$ws.Range("R8").Value = 0.130855102647
$ws.Range("S8").Value = 0.001665969177047965
$ws.Range("T8").Value = 0.001665969177047965

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gdnf"
$ws.Range("C9").Value = "Ret"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.134919666666667
$ws.Range("H9").Value = 3.404758999999999
$ws.Range("I9").Value = 0.9064795645608229
$ws.Range("J9").Value = 0.9064795645608228
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5906476666666668
$ws.Range("N9").Value = 1.771943
$ws.Range("O9").Value = 0.08473335558850956
$ws.Range("P9").Value = 0.08473335558850954
$ws.Range("Q9").Value = 0.6703376529707779
$ws.Range("R9").Value = 6.033038876737
$ws.Range("S9").Value = 0.07680905527764952
$ws.Range("T9").Value = 0.07680905527764949

